$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 802
$ws.Range("F7").Value = 101
$ws.Range("F8").Value = 10247
$ws.Range("F10").Value = 3531
$ws.Range("F11").Value = 211
$ws.Range("F12").Value = 2448
$ws.Range("F13").Value = 34
$ws.Range("F14").Value = 2817
$ws.Range("F17").Value = 2179
$ws.Range("F19").Value = 96
$ws.Range("F20").Value = 26
$ws.Range("F21").Value = 388
$ws.Range("F24").Value = 316
$ws.Range("F25").Value = 275
$ws.Range("F26").Value = 228
$ws.Range("F27").Value = 615
$ws.Range("F28").Value = 1319
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 1256
$ws.Range("F34").Value = 3711
$ws.Range("F35").Value = 3172
$ws.Range("F36").Value = 32
$ws.Range("F37").Value = 29
$ws.Range("F40").Value = 6
$ws.Range("F42").Value = 101
$ws.Range("F43").Value = 108
$ws.Range("F47").Value = 11

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 8

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 750
$ws.Range("F4").Value = 127
$ws.Range("F5").Value = 2026

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 750
$ws.Range("F4").Value = 127
$ws.Range("F9").Value = 802
$ws.Range("F10").Value = 101
$ws.Range("F11").Value = 10247
$ws.Range("F13").Value = 3531
$ws.Range("F14").Value = 211
$ws.Range("F15").Value = 2448
$ws.Range("F16").Value = 34
$ws.Range("F19").Value = 2179
$ws.Range("F21").Value = 96
$ws.Range("F22").Value = 26
$ws.Range("F23").Value = 388
$ws.Range("F25").Value = 316
$ws.Range("F26").Value = 228
$ws.Range("F27").Value = 1319
$ws.Range("F28").Value = 13
$ws.Range("F29").Value = 1256
$ws.Range("F33").Value = 8
$ws.Range("F36").Value = 3172
$ws.Range("F41").Value = 6
$ws.Range("F45").Value = 101
$ws.Range("F48").Value = 11
